# This edit inserts two new weekly price-report rows for "Zanahoria" (carrot)
# at Vega Central Mapocho de Santiago, right before the existing row 679,
# shifting all subsequent rows (old 679-705) down by two positions
# (new 681-707). The worksheet dimension is extended accordingly
# (Excel updates this automatically when rows are inserted).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at position 679-680; existing rows 679 onward
# move down to 681 onward.
$ws.Rows("679:680").Insert()

# New row 679: "Primera" quality, week of 2021-11-09 (serial 44509)
$row679 = @(
    9,
    "Vega Central Mapocho de Santiago",
    "Metropolitana",
    44509,
    13,
    100114013,
    "Zanahoria",
    "Sin especificar",
    "Primera",
    340,
    6000,
    7000,
    6500,
    "`$/saco 20 kilos",
    "Región Metropolitana",
    325,
    20,
    "Hortaliza"
)
for ($i = 0; $i -lt $row679.Length; $i++) {
    $ws.Cells.Item(679, $i + 1).Value2 = $row679[$i]
}

# New row 680: "Segunda" quality, same week (serial 44509)
$row680 = @(
    9,
    "Vega Central Mapocho de Santiago",
    "Metropolitana",
    44509,
    13,
    100114013,
    "Zanahoria",
    "Sin especificar",
    "Segunda",
    160,
    5000,
    5000,
    5000,
    "`$/saco 20 kilos",
    "Región Metropolitana",
    250,
    20,
    "Hortaliza"
)
for ($i = 0; $i -lt $row680.Length; $i++) {
    $ws.Cells.Item(680, $i + 1).Value2 = $row680[$i]
}

# Note: Rows("679:680").Insert() already copies the row-679 (now row-681)
# formatting - including the date NumberFormat on column D - down onto the
# two freshly inserted rows, so no further style fix-up is required here.
